# Layout v1 png, further development on site
#
# Table2's third column header ("Column2") is renamed to "Content", and the
# active selection moves from E6 to C2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header of the table's third column (cell C1) from "Column2" to
# "Content". Excel keeps the ListObject column name in sync with the header
# cell automatically.
$ws.Range("C1").Value = "Content"

# Update the selected cell to C2.
$ws.Range("C2").Select() | Out-Null
